$wb = $excel.ActiveWorkbook

# --- Delete the third sheet (ton-linhkien) entirely ---
$wsTon = $wb.Worksheets.Item("ton-linhkien")
$wsTon.Delete()

# --- Rename remaining sheets ---
$ws1 = $wb.Worksheets.Item("nhap-linhkien")
$ws1.Name = "nhap-thanhpham"

$ws2 = $wb.Worksheets.Item("xuat-linhkien")
$ws2.Name = "xuat-thanhpham"

# --- Clear old contents (columns A:J) on both sheets ---
$ws1.Range("A1:J10").ClearContents()
$ws2.Range("A1:J10").ClearContents()

# --- Drop the now-unused trailing columns (G:J) so the used range shrinks ---
$ws1.Columns("G:J").Delete()
$ws2.Columns("G:J").Delete()

# --- New headers (shared by both sheets) ---
$headers = @("Tên Hàng", "MCU", "Sổ Hợp Đồng", "Chip", "Ngày Nhập", "Số Lượng")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers[$i]
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data rows for nhap-thanhpham ---
# Column E holds an ISO-like date string ("2021-09-26") that must stay
# literal text rather than being auto-converted to a date serial. Force
# text format for the write, then restore the "Normal" style so no
# lingering number-format attribute is left on the cell.
$ws1.Cells.Item(2, 1).Value = "danh"
$ws1.Cells.Item(2, 2).Value = "mcu01"
$ws1.Cells.Item(2, 3).Value = "16/2021/DT-FE - 03/06/2021"
$ws1.Cells.Item(2, 4).Value = "chip01"
$ws1.Cells.Item(2, 5).NumberFormat = "@"
$ws1.Cells.Item(2, 5).Value = "2021-09-26"
$ws1.Cells.Item(2, 5).Style = "Normal"
$ws1.Cells.Item(2, 6).Value = 100

$ws1.Cells.Item(3, 1).Value = "sinh"
$ws1.Cells.Item(3, 2).Value = "mcu01"
$ws1.Cells.Item(3, 3).Value = "16/2021/DT-FE - 03/06/2021"
$ws1.Cells.Item(3, 4).Value = "chip01"
$ws1.Cells.Item(3, 5).NumberFormat = "@"
$ws1.Cells.Item(3, 5).Value = "2021-09-26"
$ws1.Cells.Item(3, 5).Style = "Normal"
$ws1.Cells.Item(3, 6).Value = 1200
